$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..T to D..U)
$ws.Columns("C").Insert()

# Header for new column C
$ws.Range("C1").Value = "custom_id_1"

# New column C values (1-based running id referencing tax-lot groupings)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 7
$ws.Range("C10").Value = 8
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 9
$ws.Range("C13").Value = 10

$ws.Range("C2").Select()
